$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update price/volume/hour cells, preserving their existing text storage type

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "278.17"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "1.14%"
$ws.Range("G2").NumberFormat = "@"
$ws.Range("G2").Value = "15"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "27.33"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "-0.51%"
$ws.Range("G3").NumberFormat = "@"
$ws.Range("G3").Value = "15"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "4.846"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "1.01%"
$ws.Range("G4").NumberFormat = "@"
$ws.Range("G4").Value = "15"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.06375"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "1.36%"
$ws.Range("G5").NumberFormat = "@"
$ws.Range("G5").Value = "15"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "7.005"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "1.17%"
$ws.Range("G6").NumberFormat = "@"
$ws.Range("G6").Value = "15"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.309"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "0.54%"
$ws.Range("G7").NumberFormat = "@"
$ws.Range("G7").Value = "15"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.8760"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-0.15%"
$ws.Range("G8").NumberFormat = "@"
$ws.Range("G8").Value = "15"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.1517"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-0.43%"
$ws.Range("G9").NumberFormat = "@"
$ws.Range("G9").Value = "15"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07506"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "0.42%"
$ws.Range("G10").NumberFormat = "@"
$ws.Range("G10").Value = "15"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.02879"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "-1.44%"
$ws.Range("G11").NumberFormat = "@"
$ws.Range("G11").Value = "15"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08954"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-0.93%"
$ws.Range("G12").NumberFormat = "@"
$ws.Range("G12").Value = "15"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.001567"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-1.37%"
$ws.Range("G13").NumberFormat = "@"
$ws.Range("G13").Value = "15"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0006390"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "0.47%"
$ws.Range("G14").NumberFormat = "@"
$ws.Range("G14").Value = "15"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.006077"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "2.88%"
$ws.Range("G15").NumberFormat = "@"
$ws.Range("G15").Value = "15"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.475"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "0.70%"
$ws.Range("G16").NumberFormat = "@"
$ws.Range("G16").Value = "15"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.297"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-0.72%"
$ws.Range("G17").NumberFormat = "@"
$ws.Range("G17").Value = "15"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.247"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "-1.10%"
$ws.Range("G18").NumberFormat = "@"
$ws.Range("G18").Value = "15"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "0.92%"
$ws.Range("G19").NumberFormat = "@"
$ws.Range("G19").Value = "15"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.05126"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "1.72%"
$ws.Range("G20").NumberFormat = "@"
$ws.Range("G20").Value = "15"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "0.44%"
$ws.Range("G21").NumberFormat = "@"
$ws.Range("G21").Value = "15"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.902"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "-0.37%"
$ws.Range("G22").NumberFormat = "@"
$ws.Range("G22").Value = "15"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.1536"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "11.30%"
$ws.Range("G23").NumberFormat = "@"
$ws.Range("G23").Value = "15"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.04400"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "0.06%"
$ws.Range("G24").NumberFormat = "@"
$ws.Range("G24").Value = "15"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.001176"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "0.20%"
$ws.Range("G25").NumberFormat = "@"
$ws.Range("G25").Value = "15"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "1.63%"
$ws.Range("G26").NumberFormat = "@"
$ws.Range("G26").Value = "15"
$ws.Range("G27").NumberFormat = "@"
$ws.Range("G27").Value = "15"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0001180"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "-1.68%"
$ws.Range("G28").NumberFormat = "@"
$ws.Range("G28").Value = "15"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "-15.08%"
$ws.Range("G29").NumberFormat = "@"
$ws.Range("G29").Value = "15"
$ws.Range("G30").NumberFormat = "@"
$ws.Range("G30").Value = "15"
$ws.Range("G31").NumberFormat = "@"
$ws.Range("G31").Value = "15"
$ws.Range("G32").NumberFormat = "@"
$ws.Range("G32").Value = "15"
$ws.Range("G33").NumberFormat = "@"
$ws.Range("G33").Value = "15"
$ws.Range("G34").NumberFormat = "@"
$ws.Range("G34").Value = "15"
$ws.Range("G35").NumberFormat = "@"
$ws.Range("G35").Value = "15"
$ws.Range("G36").NumberFormat = "@"
$ws.Range("G36").Value = "15"
$ws.Range("G37").NumberFormat = "@"
$ws.Range("G37").Value = "15"
$ws.Range("G38").NumberFormat = "@"
$ws.Range("G38").Value = "15"
$ws.Range("G39").NumberFormat = "@"
$ws.Range("G39").Value = "15"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04075"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "-0.79%"
$ws.Range("G40").NumberFormat = "@"
$ws.Range("G40").Value = "15"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.006759"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "-2.31%"
$ws.Range("G41").NumberFormat = "@"
$ws.Range("G41").Value = "15"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1403"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "20.08%"
$ws.Range("G42").NumberFormat = "@"
$ws.Range("G42").Value = "15"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-3.09%"
$ws.Range("G43").NumberFormat = "@"
$ws.Range("G43").Value = "15"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.01168"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "4.27%"
$ws.Range("G44").NumberFormat = "@"
$ws.Range("G44").Value = "15"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005351"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "3.39%"
$ws.Range("G45").NumberFormat = "@"
$ws.Range("G45").Value = "15"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.628"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "9.53%"
$ws.Range("G46").NumberFormat = "@"
$ws.Range("G46").Value = "15"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.01850"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "-7.47%"
$ws.Range("G47").NumberFormat = "@"
$ws.Range("G47").Value = "15"
$ws.Range("G48").NumberFormat = "@"
$ws.Range("G48").Value = "15"
$ws.Range("G49").NumberFormat = "@"
$ws.Range("G49").Value = "15"
$ws.Range("G50").NumberFormat = "@"
$ws.Range("G50").Value = "15"
$ws.Range("G51").NumberFormat = "@"
$ws.Range("G51").Value = "15"
